$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "Elia tomasoni "
$ws.Range("B69").Value = "Stefano Tita | Clitoriders"
$ws.Range("C69").Value = "Luca Frasca | Clitoriders"
$ws.Range("D69").Value = "Giovanni Giusto | demobusters"
$ws.Range("E69").Value = "Saverio Santoro | U.SGUARNA"
$ws.Range("F69").Value = "Moris Benedetti | Gli Introvabili"
